$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-converted to numbers by Excel (to preserve the original text layout,
# e.g. trailing zeros like "1.00" or "3.00").
$textCells = @("D5", "D7", "D8", "D10", "D17", "D18", "D19", "D21", "D23", "D25", "D28", "D32", "D34", "D37", "D38", "D40", "D42", "D44", "D45", "D46", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values.
$ws.Range("D2").Value = "28.470.97"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.555.47"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "210.88"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "24.22"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").Value = "0.0583"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "1.778.40"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "1.563.99"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "28.475.72"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "61.22"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "229.46"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").Value = "7.37"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "0.0₃0672"
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").Value = "8.92"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "150.55"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("E30").Value = "  -3.36%  "
$ws.Range("E31").Value = "  -4.32%  "
$ws.Range("D32").Value = "3.16"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").Value = "1.393.65"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "3.00"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("E35").Value = "  -4.15%  "
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.30"
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.66"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").Value = "1.94"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("D44").Value = "0.0466"
$ws.Range("E44").Value = "  +2.08%  "
$ws.Range("D45").Value = "64.48"
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("D46").Value = "5.33"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("D47").Value = "1.691.36"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("E48").Value = "  -6.56%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "85.38"
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "43.44"
$ws.Range("E50").Value = "  +5.93%  "
$ws.Range("D51").Value = "0.0₆0102"
$ws.Range("E51").Value = "  +0.35%  "
